$wb = $excel.ActiveWorkbook

# ALC row 121
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H121").Value = 2670.96
$ws.Range("I121").Value = 595
$ws.Range("J121").Value = 2757.4583
$ws.Range("K121").Value = 1785
$ws.Range("L121").Value = 8272.374899999999
$ws.Range("M121").Value = -38
$ws.Range("N121").Value = -11766.3749

# ALC row 132
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 2883.7708
$ws.Range("I132").Value = 2984.2092
$ws.Range("J132").Value = 2020
$ws.Range("K132").Value = 8952.6276
$ws.Range("L132").Value = 6060
$ws.Range("M132").Value = -6422.6276
$ws.Range("N132").Value = -11120

# ALC row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 1279.1666
$ws.Range("I138").Value = 654.5714
$ws.Range("J138").Value = 2960.7693
$ws.Range("K138").Value = 1963.7142
$ws.Range("L138").Value = 8882.3079
$ws.Range("M138").Value = 3176.2858
$ws.Range("N138").Value = -19162.3079

# ALC row 141
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 1156.1621
$ws.Range("I141").Value = 581.3570999999999
$ws.Range("J141").Value = 2944.4443
$ws.Range("K141").Value = 1744.0713
$ws.Range("L141").Value = 8833.332900000001
$ws.Range("M141").Value = 3435.9287
$ws.Range("N141").Value = -19193.3329

# ARM row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1916708.4
$ws.Range("I61").Value = 2416428
$ws.Range("K61").Value = 2416428
$ws.Range("M61").Value = -2416216

# ARM row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 8562953
$ws.Range("I132").Value = 11555044
$ws.Range("J132").Value = 85364.75
$ws.Range("K132").Value = 34665132
$ws.Range("L132").Value = 256094.25
$ws.Range("M132").Value = -34662602
$ws.Range("N132").Value = -261154.25

# ARM row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 1916708.4
$ws.Range("I136").Value = 2416428
$ws.Range("K136").Value = 7249284
$ws.Range("M136").Value = -7246734

# BSM row 80
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 290.72223
$ws.Range("J80").Value = 299.5
$ws.Range("L80").Value = 299.5
$ws.Range("N80").Value = -2295.5

# BSM row 83
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H83").Value = 290.72223
$ws.Range("J83").Value = 299.5
$ws.Range("L83").Value = 1497.5
$ws.Range("N83").Value = -11481.5

# BSM row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2958056.5
$ws.Range("I134").Value = 3297369.5
$ws.Range("J134").Value = 1185.7142
$ws.Range("K134").Value = 9892108.5
$ws.Range("L134").Value = 3557.1426
$ws.Range("M134").Value = -9889573.5
$ws.Range("N134").Value = -8627.142599999999

# CRP row 11
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H11").Value = 0
$ws.Range("J11").Value = 0
$ws.Range("L11").Value = 0
$ws.Range("N11").ClearContents()

# CRP row 58
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1333.9701
$ws.Range("I58").Value = 988.77356
$ws.Range("K58").Value = 988.77356
$ws.Range("M58").Value = -785.77356

# CRP row 122
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 952.7917
$ws.Range("I122").Value = 848.35
$ws.Range("K122").Value = 2545.05
$ws.Range("M122").Value = -95.05000000000018

# CRP row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 1805.3846
$ws.Range("I132").Value = 1362.3889
$ws.Range("J132").Value = 2802.125
$ws.Range("K132").Value = 4087.1667
$ws.Range("L132").Value = 8406.375
$ws.Range("M132").Value = -1557.1667
$ws.Range("N132").Value = -13466.375

# CRP row 134
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 1137.12
$ws.Range("I134").Value = 1254.973
$ws.Range("J134").Value = 801.6923
$ws.Range("K134").Value = 3764.919
$ws.Range("L134").Value = 2405.0769
$ws.Range("M134").Value = -1229.919
$ws.Range("N134").Value = -7475.0769

# CRP row 136
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 1333.9701
$ws.Range("I136").Value = 988.77356
$ws.Range("K136").Value = 2966.32068
$ws.Range("M136").Value = -416.3206799999998

# CUL row 5
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 8621488
$ws.Range("I5").Value = 405.4762
$ws.Range("K5").Value = 1216.4286
$ws.Range("M5").Value = -1104.4286

# CUL row 16
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H16").Value = 13114.571
$ws.Range("I16").Value = 13114.571
$ws.Range("K16").Value = 39343.713
$ws.Range("M16").Value = -39170.713

# CUL row 33
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 45454700
$ws.Range("I33").Value = 50000140
$ws.Range("J33").Value = 300
$ws.Range("K33").Value = 300000840
$ws.Range("L33").Value = 1800
$ws.Range("M33").Value = -300000557
$ws.Range("N33").Value = -2366

# CUL row 75
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H75").Value = 1005.6
$ws.Range("I75").Value = 1003.25
$ws.Range("J75").Value = 1015
$ws.Range("K75").Value = 3009.75
$ws.Range("L75").Value = 3045
$ws.Range("M75").Value = -2011.75
$ws.Range("N75").Value = -5041

# CUL row 78
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H78").Value = 1005.6
$ws.Range("I78").Value = 1003.25
$ws.Range("J78").Value = 1015
$ws.Range("K78").Value = 9029.25
$ws.Range("L78").Value = 9135
$ws.Range("M78").Value = -4037.25
$ws.Range("N78").Value = -19119

# CUL row 131
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 2693.2932
$ws.Range("J131").Value = 1856.62
$ws.Range("L131").Value = 5569.86
$ws.Range("N131").Value = -15649.86

# CUL row 135
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 8621488
$ws.Range("I135").Value = 405.4762
$ws.Range("K135").Value = 3649.2858
$ws.Range("M135").Value = -1114.2858

# GSM row 122
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 3689.925
$ws.Range("I122").Value = 4005.5186
$ws.Range("J122").Value = 3034.4614
$ws.Range("K122").Value = 12016.5558
$ws.Range("L122").Value = 9103.3842
$ws.Range("M122").Value = -9566.5558
$ws.Range("N122").Value = -14003.3842

# GSM row 126
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2067.6667
$ws.Range("I126").Value = 1901.5
$ws.Range("J126").Value = 2400
$ws.Range("K126").Value = 5704.5
$ws.Range("L126").Value = 7200
$ws.Range("M126").Value = -3234.5
$ws.Range("N126").Value = -12140

# GSM row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 18869620
$ws.Range("I132").Value = 25001504
$ws.Range("K132").Value = 75004512
$ws.Range("M132").Value = -75001982

# LTW row 122
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 1989.6364
$ws.Range("I122").Value = 1954.5
$ws.Range("J122").Value = 2083.3333
$ws.Range("K122").Value = 5863.5
$ws.Range("L122").Value = 6249.999899999999
$ws.Range("M122").Value = -3413.5
$ws.Range("N122").Value = -11149.9999

# LTW row 132
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2893.3872
$ws.Range("I132").Value = 2786.644
$ws.Range("J132").Value = 4992.6665
$ws.Range("K132").Value = 8359.931999999999
$ws.Range("L132").Value = 14977.9995
$ws.Range("M132").Value = -5829.931999999999
$ws.Range("N132").Value = -20037.9995

# WVR row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 6650121
$ws.Range("I132").Value = 7544980
$ws.Range("J132").Value = 2600.4285
$ws.Range("K132").Value = 22634940
$ws.Range("L132").Value = 7801.2855
$ws.Range("M132").Value = -22632410
$ws.Range("N132").Value = -12861.2855
